$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.5
$ws.Range("C3").Value = 7.5
$ws.Range("C3").Select()
